# Applies the "OD matrix for students" update:
#  - highlights E15/I15 in row 15 with colored fills (instead of colored font)
#  - adds a new "OD (known) assaginment 1a" section (rows 19-22) below the table
#  - tweaks sheet view (zoom + selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OD")

# --- 1. Re-style the row-15 totals for Hasselt (E15) and Leuven (I15): ---
#     switch from colored text on white to white/default text on a colored fill.
$e15 = $ws.Range("E15")
$e15.Font.ThemeColor = 1
$e15.Interior.ThemeColor = 10

$i15 = $ws.Range("I15")
$i15.Font.ThemeColor = 1
$i15.Interior.Color = 255

# --- 2. New section header (row 19), merged across B:L, bold 14pt, bordered & centered ---
$ws.Range("B19").Value = "OD (known) assaginment 1a"

$headerRow = $ws.Range("B19:L19")
$headerRow.Merge()
$headerRow.HorizontalAlignment = -4108
$headerRow.Borders.LineStyle = 1
$ws.Rows.Item(19).RowHeight = 19

$titleCell = $ws.Range("B19")
$titleCell.Font.Bold = $true
$titleCell.Font.Size = 14
$titleCell.Font.Name = "Calibri (Body)"

# --- 3. Column labels (row 20) ---
$ws.Range("C20").Value = "arriving"
$ws.Range("D20").Value = "through"

# --- 4. Leuven row (row 21) ---
$ws.Range("B21").Value = "Leuven"
$ws.Range("C21").Formula = "=14474+7718+4890"
$ws.Range("D21").Formula = "=43423+23155"
$ws.Range("C21").Interior.Color = 255

$e22src = $ws.Range("E22")
$e22src.Copy()
$ws.Range("E21").PasteSpecial(-4122)

# --- 5. Hasselt row (row 22) ---
$ws.Range("B22").Value = "Hasselt"
$ws.Range("C22").Formula = "=3491+8422"
$ws.Range("D22").Formula = "=1745+4227"
$ws.Range("C22").Interior.ThemeColor = 10

$excel.CutCopyMode = 0

# --- 6. Sheet view: zoom out and move the selection ---
$win = $excel.ActiveWindow
$win.Zoom = 67
$ws.Range("M15").Select()
